$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the Name/Role columns' contents (column formatting/width stays,
# only column A has data going forward).
$ws.Range("B1:C3").ClearContents()

# Update the three email rows with the new addresses.
$ws.Range("A2").Value = "bocaioandoru12+2@gmail.com"
$ws.Range("A3").Value = "bocaioandoru12+3@gmail.com"
$ws.Range("A4").Value = "bocaioandoru12+4@gmail.com"

# Wire up mailto hyperlinks for the two new rows (A2 already has one).
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:bocaioandoru12+3@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:bocaioandoru12+4@gmail.com")

# Keep all email cells styled as hyperlinks (Add() above also restyles A3/A4,
# this just makes sure everything lines up on the built-in Hyperlink style).
$ws.Range("A2:A4").Style = "Hyperlink"

# Match the saved selection state.
[void]$ws.Range("A3").Select()
